$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RateSheetManagement")

# Update rate-sheet staff assignments with the new roster
$ws.Range("B6").Value = "Bryce Schilling"
$ws.Range("B7").Value = "Raj Desai"
$ws.Range("B8").Value = "Gordon Bolton"
$ws.Range("B9").Value = "Blake Dickey"
$ws.Range("B4").Value = "Bryan Walker"
$ws.Range("B2").Value = "Joseph W. Swanson"

# Clear out the stale note / spare-name cells
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("D9").ClearContents()

# Make RateSheetManagement the active tab / selection, as last touched
$ws.Activate() | Out-Null
$ws.Range("C18").Select() | Out-Null
